$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the original sheet "Tabelle1" -> "1"
# ---------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item(1)
$sheet1.Name = "1"

# ---------------------------------------------------------------------
# 2. Sheet "1": two status cells flip from "in Bearbeitung" to "Erledigt"
# ---------------------------------------------------------------------
$sheet1.Range("D8").Value = "Erledigt"
$sheet1.Range("D10").Value = "Erledigt"
$sheet1.Range("D13").Value = "Erledigt"

# sheet "1" is no longer the active tab, and the remembered selection moves
$sheet1.Range("B29").Select()

# ---------------------------------------------------------------------
# 3. Duplicate sheet "1" right after itself to become sheet "2" -- this
#    carries over all styles/merges/column widths identically.
# ---------------------------------------------------------------------
$sheet1.Copy([System.Reflection.Missing]::Value, $sheet1)
$sheet2 = $wb.Worksheets.Item(2)
$sheet2.Name = "2"

# ---------------------------------------------------------------------
# 4. Sheet "2": rewrite the backlog content for Sprint Backlog 2
# ---------------------------------------------------------------------
$sheet2.Range("A1").Value = "Bimaru: Sprint Backlog 2"

$sheet2.Range("A4").Value = 5
$sheet2.Range("B4").Value = "Der Spieler kann jederzeit seine Lösung  überprüfen lassen."
$sheet2.Range("B5").ClearContents()
$sheet2.Range("C5").ClearContents()
$sheet2.Range("D5").ClearContents()
$sheet2.Range("B6").ClearContents()
$sheet2.Range("C6").ClearContents()
$sheet2.Range("D6").ClearContents()

$sheet2.Range("A7").Value = 7
$sheet2.Range("B7").Value = "Der Spieler wählt vorgegebene Spieledateien aus, um das Spiel zu starten."
$sheet2.Range("B8").ClearContents()
$sheet2.Range("C8").ClearContents()
$sheet2.Range("D8").ClearContents()
$sheet2.Range("B9").ClearContents()
$sheet2.Range("C9").ClearContents()
$sheet2.Range("D9").ClearContents()
$sheet2.Range("B10").ClearContents()
$sheet2.Range("C10").ClearContents()
$sheet2.Range("D10").ClearContents()
$sheet2.Range("B11").ClearContents()
$sheet2.Range("C11").ClearContents()
$sheet2.Range("D11").ClearContents()

$sheet2.Range("A12").Value = 2
$sheet2.Range("B12").Value = "Der Entwickler kann Spiele in XML-Dateien definieren. "
$sheet2.Range("B13").Value = "Einige Spieldateien in XML erstellen"
$sheet2.Range("C13").ClearContents()
$sheet2.Range("D13").ClearContents()
$sheet2.Range("B14").ClearContents()
$sheet2.Range("C14").ClearContents()
$sheet2.Range("D14").ClearContents()

$sheet2.Range("A15").Value = 6
$sheet2.Range("B15").Value = "Der Spieler sieht welche Schiffe zur Verfügung stehen und welche davon schon platziert wurden."
$sheet2.Range("B16").ClearContents()
$sheet2.Range("C16").ClearContents()
$sheet2.Range("D16").ClearContents()
$sheet2.Range("B17").ClearContents()
$sheet2.Range("C17").ClearContents()
$sheet2.Range("D17").ClearContents()

# rows 8, 15, 17 had custom heights copied from sheet "1"; sheet "2"'s text
# is short enough that Excel would recompute (shrink) them back to default
$sheet2.Rows.Item(8).AutoFit()
$sheet2.Rows.Item(15).AutoFit()
$sheet2.Rows.Item(17).AutoFit()

# new footer row with the Scrum Master credit, merged like the title row
$sheet2.Range("A19").Value = "Scrum Master: Remo Koller"
$sheet2.Range("A19").HorizontalAlignment = -4131
$sheet2.Range("B19").HorizontalAlignment = -4131
$sheet2.Range("A19:B19").Merge()

# sheet "2" becomes the active tab with this remembered selection
$sheet2.Range("B22").Select()
$wb.Worksheets.Item(2).Activate()
